# Logged Week 15 and simulated Week 16
# Update the "R" row (row 3) Target Depth Data on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# OFF sheet
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("B3").Value = 443
$ws.Range("C3").Value = 309
$ws.Range("D3").Value = 98
$ws.Range("E3").Value = 55

# DEF sheet
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("B3").Value = 505
$ws.Range("C3").Value = 374
$ws.Range("D3").Value = 127
$ws.Range("E3").Value = 59
